$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 448979.72
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 448979.72
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1346939.16
$ws.Range("N17").Value = -1347275.16
# Row 18
$ws.Range("H18").Value = 3277.8
$ws.Range("I18").Value = 449.5
$ws.Range("J18").Value = 5163.3335
$ws.Range("K18").Value = 449.5
$ws.Range("L18").Value = 5163.3335
$ws.Range("M18").Value = -165.5
$ws.Range("N18").Value = -5731.3335
# Row 19
$ws.Range("H19").Value = 1533.1666
$ws.Range("I19").Value = 1450
$ws.Range("J19").Value = 1574.75
$ws.Range("K19").Value = 1450
$ws.Range("L19").Value = 1574.75
$ws.Range("M19").Value = -1275
$ws.Range("N19").Value = -1924.75
# Row 20
$ws.Range("H20").Value = 2021
$ws.Range("I20").Value = 2021
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2021
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1791
$ws.Range("N20").ClearContents()
# Row 35
$ws.Range("H35").Value = 2021
$ws.Range("I35").Value = 2021
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2021
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1642
$ws.Range("N35").ClearContents()
# Row 53
$ws.Range("H53").Value = 252.77777
$ws.Range("I53").Value = 207.85715
$ws.Range("J53").Value = 410
$ws.Range("K53").Value = 207.85715
$ws.Range("L53").Value = 410
$ws.Range("M53").Value = 429.14285
$ws.Range("N53").Value = -1684
# Row 88
$ws.Range("H88").Value = 26365796
$ws.Range("I88").Value = 166668130
$ws.Range("J88").Value = 2982075.2
$ws.Range("K88").Value = 166668130
$ws.Range("L88").Value = 2982075.2
$ws.Range("M88").Value = -166667724
$ws.Range("N88").Value = -2982887.2
# Row 91
$ws.Range("H91").Value = 26365796
$ws.Range("I91").Value = 166668130
$ws.Range("J91").Value = 2982075.2
$ws.Range("K91").Value = 166668130
$ws.Range("L91").Value = 2982075.2
$ws.Range("M91").Value = -166666726
$ws.Range("N91").Value = -2984883.2
# Row 107
$ws.Range("H107").Value = 1080.9231
$ws.Range("I107").Value = 1080.9231
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1080.9231
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 839.0769
# Row 129
$ws.Range("H129").Value = 6684.875
$ws.Range("I129").Value = 1354.3572
$ws.Range("J129").Value = 43998.5
$ws.Range("K129").Value = 4063.0716
$ws.Range("L129").Value = 131995.5
$ws.Range("M129").Value = 936.9284000000002
# Row 137
$ws.Range("H137").Value = 3418.3333
$ws.Range("I137").Value = 2970.625
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 8911.875
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -6361.875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3168.2368
$ws.Range("I32").Value = 1639.1177
$ws.Range("J32").Value = 16165.75
$ws.Range("K32").Value = 1639.1177
$ws.Range("L32").Value = 16165.75
$ws.Range("M32").Value = -1352.1177
# Row 74
$ws.Range("H74").Value = 47625576
$ws.Range("I74").Value = 47625576
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 47625576
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -47624702
# Row 77
$ws.Range("H77").Value = 47625576
$ws.Range("I77").Value = 47625576
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 238127880
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -238123512
# Row 132
$ws.Range("H132").Value = 3706229
$ws.Range("I132").Value = 4002327.2
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 12006981.6
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -12004451.6
$ws.Range("N132").Value = -20055.5

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 39000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 39000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 39000
$ws.Range("N54").Value = -39968
# Row 88
$ws.Range("H88").Value = 61498.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 61498.25
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 61498.25
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -62310.25
# Row 91
$ws.Range("H91").Value = 61498.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 61498.25
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 61498.25
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -64306.25
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7586.909
$ws.Range("I31").Value = 4933.4346
$ws.Range("J31").Value = 13689.9
$ws.Range("K31").Value = 4933.4346
$ws.Range("L31").Value = 13689.9
$ws.Range("M31").Value = -4638.4346
$ws.Range("N31").Value = -14279.9
# Row 34
$ws.Range("H34").Value = 7586.909
$ws.Range("I34").Value = 4933.4346
$ws.Range("J34").Value = 13689.9
$ws.Range("K34").Value = 4933.4346
$ws.Range("L34").Value = 13689.9
$ws.Range("M34").Value = -4731.4346
$ws.Range("N34").Value = -14093.9
# Row 58
$ws.Range("H58").Value = 29419604
$ws.Range("I58").Value = 45465150
$ws.Range("J58").Value = 2764
$ws.Range("K58").Value = 45465150
$ws.Range("L58").Value = 2764
$ws.Range("M58").Value = -45464947
# Row 99
$ws.Range("H99").Value = 2994.25
$ws.Range("I99").Value = 2960.6155
$ws.Range("J99").Value = 3140
$ws.Range("K99").Value = 2960.6155
$ws.Range("L99").Value = 3140
$ws.Range("M99").Value = -1462.6155
$ws.Range("N99").Value = -6136
# Row 126
$ws.Range("H126").Value = 2994.25
$ws.Range("I126").Value = 2960.6155
$ws.Range("J126").Value = 3140
$ws.Range("K126").Value = 8881.8465
$ws.Range("L126").Value = 9420
$ws.Range("M126").Value = -6411.8465
$ws.Range("N126").Value = -14360
# Row 136
$ws.Range("H136").Value = 29419604
$ws.Range("I136").Value = 45465150
$ws.Range("J136").Value = 2764
$ws.Range("K136").Value = 136395450
$ws.Range("L136").Value = 8292
$ws.Range("M136").Value = -136392900

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 223038.75
$ws.Range("I4").Value = 250443.6
$ws.Range("J4").Value = 3800
$ws.Range("K4").Value = 751330.8
$ws.Range("L4").Value = 11400
$ws.Range("M4").Value = -751218.8
# Row 5
$ws.Range("H5").Value = 126336.625
$ws.Range("I5").Value = 167465.83
$ws.Range("J5").Value = 2949
$ws.Range("K5").Value = 502397.49
$ws.Range("L5").Value = 8847
$ws.Range("M5").Value = -502285.49
$ws.Range("N5").Value = -9071
# Row 49
$ws.Range("H49").Value = 900
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 900
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 2700
$ws.Range("N49").Value = -3012
# Row 86
$ws.Range("H86").Value = 428.6
$ws.Range("I86").Value = 373.55554
$ws.Range("J86").Value = 511.16666
$ws.Range("K86").Value = 1120.66662
$ws.Range("L86").Value = 1533.49998
$ws.Range("M86").Value = 65.33338000000003
$ws.Range("N86").Value = -3905.49998
# Row 89
$ws.Range("H89").Value = 428.6
$ws.Range("I89").Value = 373.55554
$ws.Range("J89").Value = 511.16666
$ws.Range("K89").Value = 3361.99986
$ws.Range("L89").Value = 4600.49994
$ws.Range("M89").Value = 2566.00014
$ws.Range("N89").Value = -16456.49994
# Row 131
$ws.Range("H131").Value = 1594.826
$ws.Range("I131").Value = 1165.4166
$ws.Range("J131").Value = 2063.2727
$ws.Range("K131").Value = 3496.2498
$ws.Range("L131").Value = 6189.8181
$ws.Range("M131").Value = 1543.7502
$ws.Range("N131").Value = -16269.8181
# Row 135
$ws.Range("H135").Value = 126336.625
$ws.Range("I135").Value = 167465.83
$ws.Range("J135").Value = 2949
$ws.Range("K135").Value = 1507192.47
$ws.Range("L135").Value = 26541
$ws.Range("M135").Value = -1504657.47
$ws.Range("N135").Value = -31611
# Row 137
$ws.Range("H137").Value = 8335527
$ws.Range("I137").Value = 12501590
$ws.Range("J137").Value = 3402.5
$ws.Range("K137").Value = 37504770
$ws.Range("L137").Value = 10207.5
$ws.Range("M137").Value = -37499670
$ws.Range("N137").Value = -20407.5
# Row 139
$ws.Range("H139").Value = 1411.6666
$ws.Range("I139").Value = 1412.6
$ws.Range("J139").Value = 1400
$ws.Range("K139").Value = 4237.799999999999
$ws.Range("L139").Value = 4200
$ws.Range("M139").Value = 902.2000000000007

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 70
$ws.Range("H70").Value = 9881.637000000001
$ws.Range("I70").Value = 7107.0713
$ws.Range("J70").Value = 14737.125
$ws.Range("K70").Value = 7107.0713
$ws.Range("L70").Value = 14737.125
$ws.Range("M70").Value = -6837.0713
# Row 73
$ws.Range("H73").Value = 9881.637000000001
$ws.Range("I73").Value = 7107.0713
$ws.Range("J73").Value = 14737.125
$ws.Range("K73").Value = 7107.0713
$ws.Range("L73").Value = 14737.125
$ws.Range("M73").Value = -6171.0713
# Row 107
$ws.Range("H107").Value = 1890.2106
$ws.Range("I107").Value = 1721.3334
$ws.Range("J107").Value = 2523.5
$ws.Range("K107").Value = 1721.3334
$ws.Range("L107").Value = 2523.5
$ws.Range("M107").Value = 198.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2840.25
$ws.Range("I136").Value = 2462.0908
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 7386.2724
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -4836.2724

$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108
# Row 62
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4248
# Row 65
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21240
# Row 81
$ws.Range("H81").Value = 115590
$ws.Range("I81").Value = 146902.14
$ws.Range("J81").Value = 5997.5
$ws.Range("K81").Value = 293804.28
$ws.Range("L81").Value = 11995
$ws.Range("M81").Value = -292743.28
# Row 84
$ws.Range("H84").Value = 115590
$ws.Range("I84").Value = 146902.14
$ws.Range("J84").Value = 5997.5
$ws.Range("K84").Value = 1469021.4
$ws.Range("L84").Value = 59975
$ws.Range("M84").Value = -1463717.4
# Row 125
$ws.Range("H125").Value = 119999
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 119999
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 119999
$ws.Range("N125").Value = -129839
# Row 132
$ws.Range("H132").Value = 14708935
$ws.Range("I132").Value = 17858778
$ws.Range("J132").Value = 9663.333000000001
$ws.Range("K132").Value = 53576334
$ws.Range("L132").Value = 28989.999
$ws.Range("M132").Value = -53573804
# Row 136
$ws.Range("H136").Value = 23810874
$ws.Range("I136").Value = 23810874
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 71432622
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -71430072
